$wb = $excel.ActiveWorkbook

# Sheet 1
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 134
$ws1.Range("F3").Value = 128
$ws1.Range("F4").Value = 1285
$ws1.Range("F7").Value = 987
$ws1.Range("F11").Value = 101
$ws1.Range("F12").Value = 416
$ws1.Range("F14").Value = 936
$ws1.Range("F15").Value = 1815
$ws1.Range("F16").Value = 4102
$ws1.Range("F17").Value = 1201
$ws1.Range("F18").Value = 116
$ws1.Range("F19").Value = 2664
$ws1.Range("F22").Value = 3635
$ws1.Range("F23").Value = 784
$ws1.Range("F25").Value = 45
$ws1.Range("F26").Value = 2339
$ws1.Range("F28").Value = 861
$ws1.Range("F29").Value = 172
$ws1.Range("F30").Value = 829
$ws1.Range("F31").Value = 218
$ws1.Range("F33").Value = 1376
$ws1.Range("F34").Value = 1973
$ws1.Range("F35").Value = 3
$ws1.Range("F36").Value = 501
$ws1.Range("F37").Value = 72
$ws1.Range("F39").Value = 595
$ws1.Range("F41").Value = 93
$ws1.Range("F43").Value = 238
$ws1.Range("F44").Value = 83

# Sheet 2
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 152

# Sheet 3
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 452

# Sheet 4
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 452
$ws4.Range("F3").Value = 134
$ws4.Range("F4").Value = 1285
$ws4.Range("F6").Value = 987
$ws4.Range("F8").Value = 152
$ws4.Range("F13").Value = 101
$ws4.Range("F14").Value = 416
$ws4.Range("F15").Value = 936
$ws4.Range("F16").Value = 1815
$ws4.Range("F17").Value = 4103
$ws4.Range("F18").Value = 1201
$ws4.Range("F19").Value = 116
$ws4.Range("F21").Value = 2664
$ws4.Range("F24").Value = 3635
$ws4.Range("F25").Value = 784
$ws4.Range("F28").Value = 45
$ws4.Range("F29").Value = 2339
$ws4.Range("F33").Value = 861
$ws4.Range("F34").Value = 172
$ws4.Range("F35").Value = 829
$ws4.Range("F36").Value = 218
$ws4.Range("F38").Value = 1376
$ws4.Range("F39").Value = 1973
$ws4.Range("F43").Value = 501
$ws4.Range("F44").Value = 72
$ws4.Range("F45").Value = 595
$ws4.Range("F47").Value = 93
$ws4.Range("F49").Value = 238
$ws4.Range("F50").Value = 83
